$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F7").Value = "use restrictions"
$ws.Range("F8").Value = "use restrictions"
$ws.Range("F14").Value = "93_referral_statement"
$ws.Range("F25").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F27").Value = "ppe"
$ws.Range("F28").Value = "ppe"
$ws.Range("F30").Value = "ppe"
$ws.Range("F31").Value = "ppe"
$ws.Range("F35").Value = "application instructions || env warning - species"
$ws.Range("F42").Value = "use restrictions"
$ws.Range("F43").Value = "use restrictions"
$ws.Range("F47").Value = "135_product_information"
$ws.Range("F48").Value = "application instructions"
$ws.Range("F49").Value = "use restrictions || 135_product_information"
$ws.Range("F51").Value = "mixing"
$ws.Range("F52").Value = "mixing"
$ws.Range("F55").Value = "mixing"
$ws.Range("F57").Value = "mixing"
$ws.Range("F65").Value = "use restrictions || application instructions"
$ws.Range("F67").Value = "use restrictions"
$ws.Range("F69").Value = "application instructions"
$ws.Range("F70").Value = "application instructions"
$ws.Range("F71").Value = "application instructions"
$ws.Range("F72").Value = "application instructions"
$ws.Range("F73").Value = "application instructions"
$ws.Range("F74").Value = "application instructions"
$ws.Range("F75").Value = "application instructions"
$ws.Range("F77").Value = "application instructions"
$ws.Range("F80").Value = "use restrictions"
$ws.Range("F82").Value = "use restrictions || application instructions"
$ws.Range("F83").Value = "use restrictions"
$ws.Range("F86").Value = "application instructions"
$ws.Range("F87").Value = "application instructions"
$ws.Range("F88").Value = "application instructions"
$ws.Range("F89").Value = "application instructions"
$ws.Range("F92").Value = "use restrictions || application instructions"
$ws.Range("F93").Value = "use restrictions || application instructions"
$ws.Range("F94").Value = "application instructions"
$ws.Range("F95").Value = "use restrictions || application instructions"
$ws.Range("F96").Value = "use restrictions || application instructions"
$ws.Range("F97").Value = "application instructions"
$ws.Range("F98").Value = "application instructions"
$ws.Range("F99").Value = "use restrictions || application instructions"
$ws.Range("F100").Value = "application instructions"
$ws.Range("F101").Value = "application instructions"
$ws.Range("F103").Value = "use restrictions || application instructions"
$ws.Range("F104").Value = "application instructions"
$ws.Range("F106").Value = "use restrictions || application instructions"
$ws.Range("F108").Value = "application instructions"
$ws.Range("F111").Value = "application instructions"
$ws.Range("F113").Value = "application instructions"
$ws.Range("F116").Value = "use restrictions"
$ws.Range("F119").Value = "use restrictions || application instructions"
$ws.Range("F120").Value = "application instructions"
$ws.Range("F122").Value = "use restrictions || application instructions"
$ws.Range("F127").Value = "use restrictions || application instructions"
$ws.Range("F131").Value = "application instructions"
$ws.Range("F135").Value = "application instructions"
$ws.Range("F136").Value = "application instructions"
$ws.Range("F140").Value = "use restrictions"
$ws.Range("F147").Value = "application instructions"
$ws.Range("F155").Value = "use restrictions"
$ws.Range("F162").Value = "use restrictions"
$ws.Range("F165").Value = "application instructions"
$ws.Range("F166").Value = "application instructions"
$ws.Range("F168").Value = "application instructions"
$ws.Range("F170").Value = "use restrictions"
$ws.Range("F171").Value = "application instructions"
$ws.Range("F173").Value = "application instructions"
$ws.Range("F176").Value = "use restrictions || application instructions"
$ws.Range("F177").Value = "application instructions"
$ws.Range("F178").Value = "use restrictions || application instructions"
$ws.Range("F184").Value = "use restrictions || application instructions"
$ws.Range("F199").Value = "use restrictions"
$ws.Range("F202").Value = "use restrictions"
$ws.Range("F207").Value = "use restrictions"
$ws.Range("F208").Value = "use restrictions"
$ws.Range("F219").Value = "use restrictions"
$ws.Range("F223").Value = "154_pesticide_storage"
